$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 4.6
$ws.Range("N2").Value = 6.2
$ws.Range("P2").Value = 2.84
$ws.Range("F3").Value = 5.3
$ws.Range("G3").Value = 6.4
$ws.Range("H3").Value = 1.54
$ws.Range("I3").Value = 1.66
$ws.Range("K3").Value = 5.4
$ws.Range("H5").Value = 2.16
$ws.Range("K5").Value = 4
$ws.Range("P5").Value = 1.93
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.73
$ws.Range("F7").Value = 1.97
$ws.Range("G7").Value = 2.22
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 3.8
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 4.6
$ws.Range("P7").Value = 2.28
$ws.Range("Q7").Value = 1.45
$ws.Range("F8").Value = 2.28
$ws.Range("K8").Value = 4.1
$ws.Range("F9").Value = 5.5
$ws.Range("G9").Value = 9.199999999999999
$ws.Range("H9").Value = 1.46
$ws.Range("I9").Value = 1.57
$ws.Range("J9").Value = 4.6
$ws.Range("K9").Value = 5.5
$ws.Range("P9").Value = 2.18
$ws.Range("Q9").Value = 1.49
$ws.Range("G10").Value = 2.26
$ws.Range("H10").Value = 3.05
$ws.Range("I10").Value = 3.5
$ws.Range("J10").Value = 4.1
$ws.Range("K10").Value = 4.2
$ws.Range("Q10").Value = 1.85
$ws.Range("G13").Value = 7.2
$ws.Range("H13").Value = 1.64
$ws.Range("J13").Value = 3.55
$ws.Range("F14").Value = 2.5
$ws.Range("G14").Value = 2.88
$ws.Range("I14").Value = 3.25
$ws.Range("J14").Value = 2.96
$ws.Range("G15").Value = 6.4
$ws.Range("H15").Value = 1.55
$ws.Range("I15").Value = 1.68
$ws.Range("P16").Value = 1.16
$ws.Range("Q16").Value = 1.63
$ws.Range("F17").Value = 4.4
$ws.Range("H17").Value = 1.66
$ws.Range("I17").Value = 1.83
$ws.Range("F18").Value = 1.23
$ws.Range("G18").Value = 1.28
$ws.Range("H18").Value = 11.5
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 6.6
$ws.Range("K18").Value = 8.199999999999999
$ws.Range("N18").Value = 8
$ws.Range("O18").Value = 1.11
$ws.Range("P18").Value = 3.3
$ws.Range("Q18").Value = 1.35
$ws.Range("R18").Value = 1.94
$ws.Range("S18").Value = 1.86
$ws.Range("T18").Value = 1.71
$ws.Range("U18").Value = 2.12
$ws.Range("AA18").Value = 550
$ws.Range("AB18").Value = 15
$ws.Range("AC18").Value = 19.5
$ws.Range("X18").Value = 50
$ws.Range("Y18").Value = 70
$ws.Range("AJ18").Value = 14
$ws.Range("AN18").Value = 3.55
$ws.Range("H19").Value = 2.8
$ws.Range("I19").Value = 2.84
$ws.Range("S19").Value = 2.74
$ws.Range("AB19").Value = 14.5
